$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B7"  = 5.910000000000001
    "A9"  = -21.723
    "B12" = 5.927
    "E13" = 16.579
    "D15" = -8.270999999999999
    "E16" = 16.737
    "A18" = -22.158
    "A20" = -20.557
    "E20" = 16.598
    "E24" = 16.698
    "B26" = 5.892
    "A27" = -21.188
    "B27" = 5.695000000000001
    "B29" = 5.970999999999999
    "B37" = 8.376000000000001
    "B38" = 5.763000000000001
    "D38" = -8.668999999999999
    "E39" = 16.049
    "D44" = -7.474000000000001
    "E48" = 17.403
    "B51" = 4.897
    "D51" = -8.434000000000001
    "E52" = 16.703
    "B55" = 5.705
    "E56" = 16.759
    "D57" = -8.016000000000002
    "D63" = -7.337000000000001
    "A69" = -21.404
    "B69" = 5.991000000000001
    "B70" = 5.412000000000001
    "D70" = -6.797
    "A76" = -20.396
    "A82" = -22.23
    "B83" = 5.736
    "E84" = 16.441
    "D99" = -8.138000000000002
    "E100" = 16.489
    "E101" = 16.605
    "B102" = 6.786
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
